$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("RSD_Retrofit")
$ws.Columns("M").Delete()
